$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.237.46'
$ws.Range("E2").Value = '  +0.31%  '

# Row 3
$ws.Range("D3").Value = '1.895.25'
$ws.Range("E3").Value = '  -1.15%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.05'
$ws.Range("E5").Value = '  -2.70%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5062'
$ws.Range("E7").Value = '  -3.38%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4047'
$ws.Range("E8").Value = '  -0.78%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08298'
$ws.Range("E9").Value = '  -2.62%  '

# Row 10
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.110'
$ws.Range("E10").Value = '  -1.11%  '

# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.26'
$ws.Range("E11").Value = '  -1.46%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.20'
$ws.Range("E12").Value = '  +7.78%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.422'
$ws.Range("E13").Value = '  -0.42%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.869.43'
$ws.Range("E14").Value = '  -2.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.328'
$ws.Range("E15").Value = '  -0.79%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.04%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.83'
$ws.Range("E17").Value = '  -2.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001102'
$ws.Range("E18").Value = '  -1.06%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06487'
$ws.Range("E19").Value = '  -3.00%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.40'
$ws.Range("E20").Value = '  +0.06%  '

# Row 21
$ws.Range("E21").Value = '  +0.07%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.926'
$ws.Range("E22").Value = '  -1.46%  '

# Row 23
$ws.Range("D23").Value = '30.237.26'
$ws.Range("E23").Value = '  +0.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.29'
$ws.Range("E24").Value = '  -0.38%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.188'
$ws.Range("E25").Value = '  -1.27%  '

# Row 26
$ws.Range("D26").Value = '2.106.35'
$ws.Range("E26").Value = '  -1.38%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.60'
$ws.Range("E27").Value = '  +2.21%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.81'
$ws.Range("E28").Value = '  +0.12%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.275'
$ws.Range("E29").Value = '  -5.97%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.77'
$ws.Range("E30").Value = '  -0.18%  '

# Row 31
$ws.Range("E31").Value = '  +2.44%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1042'
$ws.Range("E32").Value = '  -2.12%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.013'
$ws.Range("E33").Value = '  +0.26%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.704'
$ws.Range("E34").Value = '  +1.79%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02454'
$ws.Range("E35").Value = '  -1.27%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.319'
$ws.Range("E36").Value = '  +2.77%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06457'
$ws.Range("E37").Value = '  -2.20%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2156'
$ws.Range("E38").Value = '  -2.34%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.189'
$ws.Range("E39").Value = '  -3.11%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.607'
$ws.Range("E40").Value = '  -3.02%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6402'
$ws.Range("E41").Value = '  -2.02%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.42'
$ws.Range("E42").Value = '  -1.68%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.215'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.26'
$ws.Range("E45").Value = '  -0.28%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5975'
$ws.Range("E46").Value = '  -2.74%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.137'
$ws.Range("E47").Value = '  +2.71%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.642'
$ws.Range("E48").Value = '  -2.85%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.66'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.215'
$ws.Range("E50").Value = '  -2.33%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.74'
$ws.Range("E51").Value = '  -1.11%  '
